$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sharer column: "Text Message" (TextMessage sharer). Add header and
# mark which base item properties it accepts ("x"), matching the same
# pattern already used for the other sharers in columns B:R.
$ws.Range("S14").Value = "Text Message"

# SHKShareTypeURL block (rows 16-24)
$ws.Range("S16").Value = "x"   # item.URL
$ws.Range("S20").Value = "x"   # item.image
$ws.Range("S21").Value = "x"   # item.title
$ws.Range("S22").Value = "x"   # item.text
$ws.Range("S24").Value = "x"   # item.data (file)

# SHKShareTypeText block (rows 30-36)
$ws.Range("S30").Value = "x"   # item.URL
$ws.Range("S32").Value = "x"   # item.image
$ws.Range("S33").Value = "x"   # item.title
$ws.Range("S34").Value = "x"   # item.text
$ws.Range("S36").Value = "x"   # item.data (file)

# SHKShareTypeImage block (rows 42-48)
$ws.Range("S42").Value = "x"   # item.URL
$ws.Range("S44").Value = "x"   # item.image
$ws.Range("S45").Value = "x"   # item.title
$ws.Range("S46").Value = "x"   # item.text
$ws.Range("S48").Value = "x"   # item.data (file)

# SHKShareTypeFile block (rows 55-61)
$ws.Range("S55").Value = "x"   # item.URL
$ws.Range("S57").Value = "x"   # item.image
$ws.Range("S58").Value = "x"   # item.title
$ws.Range("S59").Value = "x"   # item.text
$ws.Range("S61").Value = "x"   # item.file (file type)

# "custom values" separator rows have no marker for this sharer - remove the
# placeholder formatted-but-empty cells there entirely.
$ws.Range("S37").Clear()
$ws.Range("S38").Clear()
$ws.Range("S49").Clear()
$ws.Range("S50").Clear()
$ws.Range("S62").Clear()
$ws.Range("S63").Clear()

# SHKShareTypeUserInfo summary row: Text Message sharer does not support it.
$ws.Range("S67").Value = "N/A"

# Reflect where the author was last working when the sheet was saved.
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("S57:S61").Select()
